# Apply the "Upload new version with timestamp" edit:
#  1. Insert a new product row ("LAXEOL PI 5MG  250TAB") in its correct alphabetical
#     position (row 25), shifting the rest of the table (and the totals/footer rows) down
#     by one row.
#  2. Populate the new row with its balance / order-limit / price / sale-price / transaction
#     values and copy the formatting (styles + row height + cell merges) from the row above.
#  3. Refresh the running index numbers (column A) for every shifted data row plus the new
#     last row so they remain sequential (1..38).
#  4. Update the grand-total (sum of the "sale price" column) to include the new row.
#  5. Update the report-generation timestamp shown in the footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at row 25 (this pushes rows 25-45 down to 26-46) ---------------
$ws.Rows("25:25").Insert()

# Copy the formatting (number formats, fonts, fills, borders, alignment) of the row above
# (row 24, "IVYWONDER ...") onto the freshly inserted, still-blank row 25.
$ws.Range("A24:Q24").Copy()
$ws.Range("A25:Q25").PasteSpecial(-4122)

# Match the row height used by the rest of the product rows with this row's own original
# height (24.75pt, same as before the insert).
$ws.Rows("25:25").RowHeight = 24.75

# Re-create the per-row cell merges (A:B, C:G, H:K, L:M, N:O) for the new row, matching the
# merge pattern used by every other product row.
$ws.Range("A25:B25").Merge()
$ws.Range("C25:G25").Merge()
$ws.Range("H25:K25").Merge()
$ws.Range("L25:M25").Merge()
$ws.Range("N25:O25").Merge()

# --- 2. Fill in the values for the new product row --------------------------------------
$ws.Range("A25").Value2 = 19
$ws.Range("C25").Value2 = "LAXEOL PI 5MG  250TAB"
$ws.Range("H25").Value2 = "0:3"
$ws.Range("L25").Value2 = "0"
$ws.Range("N25").Value2 = "300.00"
$ws.Range("P25").Value2 = "12.0000"
$ws.Range("Q25").Value2 = "0:1"

# --- 3. Renumber column A (the running "م" index) for every row from 26 down to the new
#        last data row (44) so the sequence stays 20, 21, 22 ... 38 -----------------------
$ws.Range("A26").Value2 = 20
$ws.Range("A27").Value2 = 21
$ws.Range("A28").Value2 = 22
$ws.Range("A29").Value2 = 23
$ws.Range("A30").Value2 = 24
$ws.Range("A31").Value2 = 25
$ws.Range("A32").Value2 = 26
$ws.Range("A33").Value2 = 27
$ws.Range("A34").Value2 = 28
$ws.Range("A35").Value2 = 29
$ws.Range("A36").Value2 = 30
$ws.Range("A37").Value2 = 31
$ws.Range("A38").Value2 = 32
$ws.Range("A39").Value2 = 33
$ws.Range("A40").Value2 = 34
$ws.Range("A41").Value2 = 35
$ws.Range("A42").Value2 = 36
$ws.Range("A43").Value2 = 37
$ws.Range("A44").Value2 = 38

# --- 4. Update the grand total (column P) to include the new row's sale price -----------
$ws.Range("P45").Value2 = 1524.745

# --- 5. Update the report-generation timestamp shown in the footer ----------------------
$ws.Range("A46").Value2 = "Wednesday, 6 August, 2025 1:07 PM"
